$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.862.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +6.40%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.016.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.75%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'584.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.52%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  +12.01%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "'LidoStakedEther"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'3.012.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.80%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "'XRP"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.75%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'6.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.88%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +6.79%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +7.54%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +9.31%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'34.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.50%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'65.830.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +6.38%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'3.514.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.71%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'6.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +7.61%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.012.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.05%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'457.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.57%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'14.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +7.95%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +5.87%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'7.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +7.81%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'82.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.64%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +12.81%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'12.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.92%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'10.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.22%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  -0.04%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +15.17%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'2.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +16.22%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.0000107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.80%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  +3.74%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'27.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +6.18%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.37%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.11%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.92%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +8.74%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +14.17%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("B39").Value = "'dogwifhat"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.68%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = "'OKB"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'49.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.00%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +16.60%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +6.76%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'43.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.57%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'8.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.80%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'391.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +12.55%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.800.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.63%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.0355"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.35%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'134.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.00%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -0.07%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'23.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +9.86%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  +4.50%  "
$ws.Range("E51").Style = "Normal"
